$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need their cell
# pre-formatted as Text so Excel stores them as strings (matching the
# source data, which is always text) instead of auto-converting to a number.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.761.50'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '3.842.00'
$ws.Range("E3").Value = '  +2.72%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '609.91'
$ws.Range("E5").Value = '  -2.20%  '
$ws.Range("D6").Value = '174.47'
$ws.Range("E6").Value = '  -3.32%  '
$ws.Range("D7").Value = '3.840.22'
$ws.Range("E7").Value = '  +2.82%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.526'
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("E10").Value = '  -1.33%  '
$ws.Range("D11").Value = '6.45'
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("D12").Value = '0.479'
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("D13").Value = '39.87'
$ws.Range("E13").Value = '  -2.62%  '
$ws.Range("D14").Value = '0.0000254'
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").Value = '4.473.39'
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("D16").Value = '3.836.49'
$ws.Range("E16").Value = '  +2.65%  '
$ws.Range("D17").Value = '69.817.06'
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = '7.45'
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("E19").Value = '  -3.25%  '
$ws.Range("D20").Value = '16.60'
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").Value = '504.55'
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("D23").Value = '0.735'
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("E24").Value = '  -5.50%  '
$ws.Range("D25").Value = '85.77'
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("E26").Value = '  +4.91%  '
$ws.Range("D27").Value = '12.61'
$ws.Range("E27").Value = '  -4.21%  '
$ws.Range("D28").Value = '10.38'
$ws.Range("E28").Value = '  -9.41%  '
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("D32").Value = '7.95'
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("D33").Value = '32.26'
$ws.Range("E33").Value = '  +3.25%  '
$ws.Range("E34").Value = '  -2.62%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("E37").Value = '  -2.02%  '
$ws.Range("E38").Value = '  +2.11%  '
$ws.Range("D39").Value = '487.28'
$ws.Range("E39").Value = '  +12.97%  '
$ws.Range("E40").Value = '  -1.26%  '
$ws.Range("E41").Value = '  -2.81%  '
$ws.Range("D42").Value = '49.72'
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("E43").Value = '  +3.85%  '
$ws.Range("D44").Value = '43.11'
$ws.Range("E44").Value = '  -6.11%  '
$ws.Range("D45").Value = '8.51'
$ws.Range("E45").Value = '  -2.75%  '
$ws.Range("D46").Value = '2.924.43'
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("D47").Value = '0.0360'
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("D48").Value = '139.95'
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = '26.75'
$ws.Range("E50").Value = '  -3.12%  '
$ws.Range("D51").Value = '2.41'
$ws.Range("E51").Value = '  -4.47%  '
